$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B2 = New-Object 'object[,]' 24,3
$arr_B2[0,0] = 0.6615872312322892
$arr_B2[0,1] = 0.188700032303899
$arr_B2[0,2] = 0.02282370822750579
$arr_B2[1,0] = 0.5775129746917571
$arr_B2[1,1] = 0.1698607710714839
$arr_B2[1,2] = 0.01989670184102721
$arr_B2[2,0] = 0.5256859433023635
$arr_B2[2,1] = 0.1582286510992787
$arr_B2[2,2] = 0.0180940746111915
$arr_B2[3,0] = 0.5045161983776723
$arr_B2[3,1] = 0.1534725908710186
$arr_B2[3,2] = 0.0173581735593018
$arr_B2[4,0] = 0.5009980219919044
$arr_B2[4,1] = 0.1526819035991025
$arr_B2[4,2] = 0.0172358997609976
$arr_B2[5,0] = 0.5254006400535332
$arr_B2[5,1] = 0.1581645729611125
$arr_B2[5,2] = 0.01808415525296425
$arr_B2[6,0] = 0.6326420283974414
$arr_B2[6,1] = 0.1822179266479509
$arr_B2[6,2] = 0.02181563358323046
$arr_B2[7,0] = 0.8412468276662821
$arr_B2[7,1] = 0.2288572764969103
$arr_B2[7,2] = 0.02908805310548956
$arr_B2[8,0] = 0.9933982818477034
$arr_B2[8,1] = 0.2627828106576828
$arr_B2[8,2] = 0.03440162092563526
$arr_B2[9,0] = 1.062359323261319
$arr_B2[9,1] = 0.2781389207150085
$arr_B2[9,2] = 0.03681209182119005
$arr_B2[10,0] = 1.088435018808752
$arr_B2[10,1] = 0.2839424681983189
$arr_B2[10,2] = 0.03772386537275452
$arr_B2[11,0] = 1.082820887443688
$arr_B2[11,1] = 0.2826930878401015
$arr_B2[11,2] = 0.03752754486339427
$arr_B2[12,0] = 1.064505366044727
$arr_B2[12,1] = 0.2786166145941991
$arr_B2[12,2] = 0.03688712473451972
$arr_B2[13,0] = 1.053281528789569
$arr_B2[13,1] = 0.2761181492976732
$arr_B2[13,2] = 0.03649471447552344
$arr_B2[14,0] = 0.9888862525637023
$arr_B2[14,1] = 0.2617776702963681
$arr_B2[14,2] = 0.03424395139386149
$arr_B2[15,0] = 0.9493155673989122
$arr_B2[15,1] = 0.2529602734230707
$arr_B2[15,2] = 0.03286142564091676
$arr_B2[16,0] = 0.926531817040825
$arr_B2[16,1] = 0.2478815385548501
$arr_B2[16,2] = 0.03206560634355071
$arr_B2[17,0] = 0.9188136118387433
$arr_B2[17,1] = 0.2461607426611181
$arr_B2[17,2] = 0.03179604962157612
$arr_B2[18,0] = 0.9535304051382241
$arr_B2[18,1] = 0.2538996488564464
$arr_B2[18,2] = 0.03300866317216844
$arr_B2[19,0] = 1.069886134799731
$arr_B2[19,1] = 0.279814288588625
$arr_B2[19,2] = 0.03707525976754766
$arr_B2[20,0] = 1.145706888182247
$arr_B2[20,1] = 0.2966838986110645
$arr_B2[20,2] = 0.03972704729879695
$arr_B2[21,0] = 1.105261092587284
$arr_B2[21,1] = 0.2876865543257452
$arr_B2[21,2] = 0.03831230313872425
$arr_B2[22,0] = 0.9516249828677701
$arr_B2[22,1] = 0.2534749867685093
$arr_B2[22,2] = 0.03294210015404531
$arr_B2[23,0] = 0.7850031685738372
$arr_B2[23,1] = 0.2162985931587116
$arr_B2[23,2] = 0.02712569950052313
$ws.Range("B2:D25").Value = $arr_B2

$arr_F2 = New-Object 'object[,]' 24,4
$arr_F2[0,0] = 0.2863447698869095
$arr_F2[0,1] = 0.1502701152003141
$arr_F2[0,2] = 0.3197164537272386
$arr_F2[0,3] = 0.2047964974952876
$arr_F2[1,0] = 0.2835585504525113
$arr_F2[1,1] = 0.1486296955100741
$arr_F2[1,2] = 0.3220292755788705
$arr_F2[1,3] = 0.208827885884288
$arr_F2[2,0] = 0.2820944199426165
$arr_F2[2,1] = 0.1477947510039002
$arr_F2[2,2] = 0.3236329978720178
$arr_F2[2,3] = 0.2114812982298702
$arr_F2[3,0] = 0.2815597426138652
$arr_F2[3,1] = 0.1474976979942966
$arr_F2[3,2] = 0.3243327105834197
$arr_F2[3,3] = 0.2126072918848223
$arr_F2[4,0] = 0.2814747011655143
$arr_F2[4,1] = 0.1474509776035191
$arr_F2[4,2] = 0.3244516868229539
$arr_F2[4,3] = 0.2127969599684434
$arr_F2[5,0] = 0.2820869582545313
$arr_F2[5,1] = 0.1477905701231208
$arr_F2[5,2] = 0.3236422474430611
$arr_F2[5,3] = 0.2114963028701933
$arr_F2[6,0] = 0.2853328658738334
$arr_F2[6,1] = 0.1496686643855369
$arr_F2[6,2] = 0.3204758029657384
$arr_F2[6,3] = 0.2061495006915326
$arr_F2[7,0] = 0.293657839131356
$arr_F2[7,1] = 0.1547247636158673
$arr_F2[7,2] = 0.3157234491850076
$arr_F2[7,3] = 0.1970813996644072
$arr_F2[8,0] = 0.3009741935003376
$arr_F2[8,1] = 0.1592858169774658
$arr_F2[8,2] = 0.3131202719518029
$arr_F2[8,3] = 0.1912876775211103
$arr_F2[9,0] = 0.3045643417628128
$arr_F2[9,1] = 0.1615465268500031
$arr_F2[9,2] = 0.3121289647803067
$arr_F2[9,3] = 0.188841585576963
$arr_F2[10,0] = 0.3059615648615974
$arr_F2[10,1] = 0.1624294691628592
$arr_F2[10,2] = 0.3117813207257498
$arr_F2[10,3] = 0.1879426520683509
$arr_F2[11,0] = 0.3056589698373742
$arr_F2[11,1] = 0.1622381150015926
$arr_F2[11,2] = 0.311854958165668
$arr_F2[11,3] = 0.1881350355964653
$arr_F2[12,0] = 0.3046785360630366
$arr_F2[12,1] = 0.1616186279617366
$arr_F2[12,2] = 0.3120998078679094
$arr_F2[12,3] = 0.1887670811121609
$arr_F2[13,0] = 0.3040829047342157
$arr_F2[13,1] = 0.1612426762299606
$arr_F2[13,2] = 0.3122533984074707
$arr_F2[13,3] = 0.189157791719607
$arr_F2[14,0] = 0.3007448433851678
$arr_F2[14,1] = 0.1591418257863069
$arr_F2[14,2] = 0.3131889381018382
$arr_F2[14,3] = 0.1914513570581846
$arr_F2[15,0] = 0.2987641669610426
$arr_F2[15,1] = 0.1579007219020099
$arr_F2[15,2] = 0.3138122691027121
$arr_F2[15,3] = 0.1929069911525882
$arr_F2[16,0] = 0.2976495828221459
$arr_F2[16,1] = 0.1572043582929652
$arr_F2[16,2] = 0.3141889481701554
$arr_F2[16,3] = 0.1937620643096984
$arr_F2[17,0] = 0.297276435618528
$arr_F2[17,1] = 0.1569715807657843
$arr_F2[17,2] = 0.3143196033093716
$arr_F2[17,3] = 0.1940546368584357
$arr_F2[18,0] = 0.2989724619415384
$arr_F2[18,1] = 0.158031028766743
$arr_F2[18,2] = 0.3137440353844241
$arr_F2[18,3] = 0.1927501902697362
$arr_F2[19,0] = 0.3049654893457756
$arr_F2[19,1] = 0.1617998562795648
$arr_F2[19,2] = 0.3120271366156189
$arr_F2[19,3] = 0.188580691014522
$arr_F2[20,0] = 0.3091021148729993
$arr_F2[20,1] = 0.1644196273738032
$arr_F2[20,2] = 0.3110667503033326
$arr_F2[20,3] = 0.1860151348963921
$arr_F2[21,0] = 0.3068741887903741
$arr_F2[21,1] = 0.1630070314580507
$arr_F2[21,2] = 0.3115645288870965
$arr_F2[21,3] = 0.1873697982358653
$arr_F2[22,0] = 0.2988782166130264
$arr_F2[22,1] = 0.1579720635810702
$arr_F2[22,2] = 0.3137748267904428
$arr_F2[22,3] = 0.1928210232350143
$arr_F2[23,0] = 0.2911953554688083
$arr_F2[23,1] = 0.1532089888308619
$arr_F2[23,2] = 0.3168530867345538
$arr_F2[23,3] = 0.199382368038064
$ws.Range("F2:I25").Value = $arr_F2

$arr_N2 = New-Object 'object[,]' 24,2
$arr_N2[0,0] = 0.8159831952810066
$arr_N2[0,1] = 0.8388828661067009
$arr_N2[1,0] = 0.804681593931349
$arr_N2[1,1] = 0.8400333483516533
$arr_N2[2,0] = 0.7980274380249739
$arr_N2[2,1] = 0.8414819075362772
$arr_N2[3,0] = 0.7953879915046542
$arr_N2[3,1] = 0.8422585640705904
$arr_N2[4,0] = 0.794954087231531
$arr_N2[4,1] = 0.8423987749908832
$arr_N2[5,0] = 0.7979915486323108
$arr_N2[5,1] = 0.8414916275979039
$arr_N2[6,0] = 0.8120275249493716
$arr_N2[6,1] = 0.8391253561090082
$arr_N2[7,0] = 0.8417945568806431
$arr_N2[7,1] = 0.8403870249187122
$arr_N2[8,0] = 0.8650086772045285
$arr_N2[8,1] = 0.8449320818470341
$arr_N2[9,0] = 0.8758567141804292
$arr_N2[9,1] = 0.8477897700012988
$arr_N2[10,0] = 0.8800055340223167
$arr_N2[10,1] = 0.8489858306335094
$arr_N2[11,0] = 0.8791101998937734
$arr_N2[11,1] = 0.8487231664324639
$arr_N2[12,0] = 0.8761972230713297
$arr_N2[12,1] = 0.8478858858903777
$arr_N2[13,0] = 0.8744182532360867
$arr_N2[13,1] = 0.8473878713882783
$arr_N2[14,0] = 0.8643054826042373
$arr_N2[14,1] = 0.8447612500049928
$arr_N2[15,0] = 0.8581749963170751
$arr_N2[15,1] = 0.8433524848679355
$arr_N2[16,0] = 0.8546760221011738
$arr_N2[16,1] = 0.8426165526518474
$arr_N2[17,0] = 0.8534960021557367
$arr_N2[17,1] = 0.8423801395560986
$arr_N2[18,0] = 0.8588247938183571
$arr_N2[18,1] = 0.8434947529715799
$arr_N2[19,0] = 0.8770517286964719
$arr_N2[19,1] = 0.8481287214794264
$arr_N2[20,0] = 0.8892022808718139
$arr_N2[20,1] = 0.851821413065835
$arr_N2[21,0] = 0.8826956662326779
$arr_N2[21,1] = 0.8497896893735515
$arr_N2[22,0] = 0.8585309408261281
$arr_N2[22,1] = 0.8434302031271272
$arr_N2[23,0] = 0.8335040865935071
$arr_N2[23,1] = 0.8394117711048068
$ws.Range("N2:O25").Value = $arr_N2

Write-Output "done"